# Adjusting CPRA gate closure triggers based on feedback from CPRA.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("C11").Value = 2.5

# Leave the active selection on D1, matching the saved view state.
$ws.Range("D1").Select()
